# Applies the "Added changes for Scheduler" update:
#  - Sheet1: refreshed Camarilla pivot levels (columns G-K, rows 3-14)
#    for Copper/Zinc/Lead/Nickel/Aluminium.
#  - Sheet2: refreshed OHLC scheduler source rows for Copper/Zinc/Lead/
#    Nickel/Aluminium (rows 6-10), including rolling the expiry date
#    shown in column B from 30APR2021 to 31MAR2021.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 ("Camarilla" pivot levels) - columns G,H,I,J,K for rows 3..14
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$sheet1Updates = @(
    @{ Row = 3;  G = 688.04;            H = 219.78; I = 163.2;  J = 1177.76; K = 173.11 },
    @{ Row = 4;  G = 686.09;            H = 218.72; I = 162.52; J = 1172.78; K = 172.58 },
    @{ Row = 5;  G = 684.15;            H = 217.66; I = 161.84; J = 1167.8;  K = 172.04 },
    @{ Row = 6;  G = 681.8200000000001; H = 216.38; I = 161.05; J = 1161.75; K = 171.4  },
    @{ Row = 7;  G = 681.05;            H = 215.95; I = 160.78; J = 1159.73; K = 171.18 },
    @{ Row = 8;  G = 680.27;            H = 215.53; I = 160.52; J = 1157.72; K = 170.97 },
    @{ Row = 9;  G = 678.73;            H = 214.67; I = 159.98; J = 1153.68; K = 170.53 },
    @{ Row = 10; G = 677.95;            H = 214.25; I = 159.72; J = 1151.67; K = 170.32 },
    @{ Row = 11; G = 677.1799999999999; H = 213.82; I = 159.45; J = 1149.65; K = 170.1  },
    @{ Row = 12; G = 674.85;            H = 212.54; I = 158.66; J = 1143.6;  K = 169.46 },
    @{ Row = 13; G = 672.91;            H = 211.48; I = 157.98; J = 1138.62; K = 168.92 },
    @{ Row = 14; G = 670.96;            H = 210.42; I = 157.3;  J = 1133.64; K = 168.39 }
)

foreach ($u in $sheet1Updates) {
    $r = $u.Row
    $ws1.Range("G$r").Value = $u.G
    $ws1.Range("H$r").Value = $u.H
    $ws1.Range("I$r").Value = $u.I
    $ws1.Range("J$r").Value = $u.J
    $ws1.Range("K$r").Value = $u.K
}

# ---------------------------------------------------------------------
# Sheet2 (scheduler source data) - rows 6..10, columns B (expiry date)
# and C..G (Open/High/Low/Close/Previous Close)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$sheet2Updates = @(
    @{ Row = 6;  Date = "31MAR2021"; C = 679.8; D = 680.65; E = 672.2;  F = 679.5;  G = 684.1  },
    @{ Row = 7;  Date = "31MAR2021"; C = 218.35; D = 218.35; E = 213.7; F = 215.1;  G = 217.8  },
    @{ Row = 8;  Date = "31MAR2021"; C = 158.6; D = 160.5;  E = 157.6;  F = 160.25; G = 158.7  },
    @{ Row = 9;  Date = "31MAR2021"; C = 1174;   D = 1174.8; E = 1152.8; F = 1155.7; G = 1181.3 },
    @{ Row = 10; Date = "31MAR2021"; C = 172.4;  D = 172.4;  E = 170.05; F = 170.75; G = 172.9  }
)

foreach ($u in $sheet2Updates) {
    $r = $u.Row
    $ws2.Range("B$r").Value = $u.Date
    $ws2.Range("C$r").Value = $u.C
    $ws2.Range("D$r").Value = $u.D
    $ws2.Range("E$r").Value = $u.E
    $ws2.Range("F$r").Value = $u.F
    $ws2.Range("G$r").Value = $u.G
}
